$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.368.86"
$ws.Range("E2").Value = "  -3.05%  "

$ws.Range("D3").Value = "2.224.66"
$ws.Range("E3").Value = "  -2.15%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "109.96"
$ws.Range("E5").Value = "  -8.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "296.54"
$ws.Range("E6").Value = "  +11.50%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.619"
$ws.Range("E7").Value = "  -4.25%  "

$ws.Range("E8").Value = "  -0.39%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.602"
$ws.Range("E9").Value = "  -3.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.00"
$ws.Range("E10").Value = "  -7.99%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0915"
$ws.Range("E11").Value = "  -3.50%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.25"
$ws.Range("E12").Value = "  -0.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.81"
$ws.Range("E13").Value = "  -4.46%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.01"
$ws.Range("E14").Value = "  +11.73%  "

$ws.Range("E15").Value = "  -2.62%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.01"
$ws.Range("E16").Value = "  -3.12%  "

$ws.Range("D17").Value = "2.556.32"
$ws.Range("E17").Value = "  -2.28%  "

$ws.Range("D18").Value = "2.227.89"
$ws.Range("E18").Value = "  -2.17%  "

$ws.Range("D19").Value = "42.368.73"
$ws.Range("E19").Value = "  -2.95%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.43"
$ws.Range("E20").Value = "  +7.52%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000106"
$ws.Range("E21").Value = "  -4.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.47"
$ws.Range("E22").Value = "  +0.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.49"
$ws.Range("E23").Value = "  +21.21%  "

$ws.Range("E24").Value = "  -3.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "228.95"
$ws.Range("E25").Value = "  -3.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.31"
$ws.Range("E26").Value = "  -2.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.69"
$ws.Range("E27").Value = "  -2.55%  "

$ws.Range("E28").Value = "  -1.67%  "

$ws.Range("E29").Value = "  -0.74%  "

$ws.Range("E30").Value = "  -8.79%  "

$ws.Range("E31").Value = "  -5.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "173.88"
$ws.Range("E32").Value = "  +1.08%  "

$ws.Range("E33").Value = "  -2.85%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0893"
$ws.Range("E34").Value = "  -2.81%  "

$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.20"
$ws.Range("E35").Value = "  +13.59%  "

$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.71"
$ws.Range("E36").Value = "  -1.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.33"
$ws.Range("E37").Value = "  +2.42%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0381"
$ws.Range("E38").Value = "  -1.84%  "

$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.126"
$ws.Range("E39").Value = "  -3.86%  "

$ws.Range("E40").Value = "  -0.73%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.43"
$ws.Range("E41").Value = "  -4.97%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.12"
$ws.Range("E42").Value = "  -2.39%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.237"
$ws.Range("E43").Value = "  -0.26%  "

$ws.Range("E44").Value = "  -0.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.67"
$ws.Range("E45").Value = "  -8.41%  "

$ws.Range("E46").Value = "  -4.47%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.44"
$ws.Range("E47").Value = "  -6.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.32"
$ws.Range("E48").Value = "  +3.91%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.52"
$ws.Range("E49").Value = "  +1.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.46"
$ws.Range("E50").Value = "  -1.01%  "

$ws.Range("E51").Value = "  +5.08%  "

